$wb = $excel.ActiveWorkbook

# --- Summary sheet: widen the selected/highlighted range from row 12 to row 14 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

# --- Repayment schedule sheet: add column O values (between N and P) for rows 2-14 ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Copy the formatting of column N onto the new column O so the new
# cells pick up the same style (s="7") used throughout the sheet.
$ws.Range("N2:N14").Copy()
$ws.Range("O2:O14").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 4; $r -le 14; $r++) {
    $ws.Cells.Item($r, 15).Value = 0
}

# --- Transactions sheet: update A2:A4 values and selection ---
$ws2 = $wb.Worksheets.Item("Transactions")
$ws2.Cells.Item(2, 1).Value = 100
$ws2.Cells.Item(3, 1).Value = 98
$ws2.Cells.Item(4, 1).Value = 96
$ws2.Activate()
$ws2.Range("D4").Select()
